$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5172
$wsExhibit.Range("F5").Value = 7467
$wsExhibit.Range("F11").Value = 30
$wsExhibit.Range("F12").Value = 4324
$wsExhibit.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/k5ZRuF6i1721816401868.jpeg"
$wsExhibit.Range("F13").Value = 1765
$wsExhibit.Range("F14").Value = 105
$wsExhibit.Range("F16").Value = 2923
$wsExhibit.Range("F18").Value = 567
$wsExhibit.Range("F19").Value = 209
$wsExhibit.Range("F20").Value = 504
$wsExhibit.Range("F21").Value = 441
$wsExhibit.Range("F22").Value = 458
$wsExhibit.Range("F23").Value = 309
$wsExhibit.Range("F24").Value = 103
$wsExhibit.Range("F25").Value = 1695
$wsExhibit.Range("F26").Value = 1191
$wsExhibit.Range("F28").Value = 1384
$wsExhibit.Range("F29").Value = 108
$wsExhibit.Range("F30").Value = 580
$wsExhibit.Range("F32").Value = 516
$wsExhibit.Range("F36").Value = 67
$wsExhibit.Range("F37").Value = 2915
$wsExhibit.Range("F39").Value = 23
$wsExhibit.Range("F40").Value = 85
$wsExhibit.Range("F42").Value = 41

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 13

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5172
$wsAll.Range("F5").Value = 7467
$wsAll.Range("F11").Value = 30
$wsAll.Range("F12").Value = 4324
$wsAll.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202407/k5ZRuF6i1721816401868.jpeg"
$wsAll.Range("F13").Value = 1765
$wsAll.Range("F14").Value = 105
$wsAll.Range("F16").Value = 2923
$wsAll.Range("F18").Value = 567
$wsAll.Range("F19").Value = 209
$wsAll.Range("F20").Value = 504
$wsAll.Range("F21").Value = 441
$wsAll.Range("F22").Value = 458
$wsAll.Range("F24").Value = 309
$wsAll.Range("F25").Value = 103
$wsAll.Range("F26").Value = 1695
$wsAll.Range("F27").Value = 1191
$wsAll.Range("F29").Value = 1384
$wsAll.Range("F30").Value = 108
$wsAll.Range("F31").Value = 580
$wsAll.Range("F33").Value = 516
$wsAll.Range("F37").Value = 67
$wsAll.Range("F38").Value = 2915
$wsAll.Range("F39").Value = 13
$wsAll.Range("F41").Value = 23
$wsAll.Range("F42").Value = 85
$wsAll.Range("F44").Value = 42
